$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "floor" sheet: the domain no longer carries a "0" floor code, so
# that row is removed entirely (everything below shifts up).
# ------------------------------------------------------------------
$wsFloor = $wb.Worksheets.Item("floor")
$zeroCell = $wsFloor.Columns.Item(1).Find(0)
$zeroCell.EntireRow.Delete()

# ------------------------------------------------------------------
# "functional category" sheet: "Assignable"/"Non-Assignable" are
# renamed to "Leasable"/"Non-Leasable", and a new "Vacant" category
# is inserted right under the header.
# ------------------------------------------------------------------
$wsFuncCat = $wb.Worksheets.Item("functional category")

$wsFuncCat.Rows.Item(2).Insert()
$vacantCell = $wsFuncCat.Cells.Item(2, 1)
$vacantCell.Value = "Vacant"
$vacantCell.Font.Bold = $false

$wsFuncCat.Columns.Item(1).Find("Assignable").Value = "Leasable"
$wsFuncCat.Columns.Item(1).Find("Non-Assignable").Value = "Non-Leasable"

# Page was set up for portrait printing
$wsFuncCat.PageSetup.Orientation = 1

# ------------------------------------------------------------------
# "functional use" sheet: "Vacant" is no longer a functional use, so
# that row is removed entirely.
# ------------------------------------------------------------------
$wsFuncUse = $wb.Worksheets.Item("functional use")
$wsFuncUse.Columns.Item(1).Find("Vacant").EntireRow.Delete()
